$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells stay text (avoid Excel auto-numeric conversion
# dropping significant trailing zeros or multi-dot thousands separators).
$priceCells = @('D2', 'D3', 'D5', 'D8', 'D10', 'D11', 'D12', 'D13', 'D14', 'D15', 'D16', 'D17', 'D18', 'D19', 'D20', 'D22', 'D25', 'D26', 'D27', 'D30', 'D33', 'D34', 'D35', 'D37', 'D38', 'D39', 'D40', 'D42', 'D43', 'D44', 'D45', 'D46', 'D47', 'D48', 'D49', 'D50', 'D51')
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '34.625.82'
$ws.Range('E2').Value = '  +2.47%  '
$ws.Range('D3').Value = '1.788.17'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '223.32'
$ws.Range('E5').Value = '  -0.96%  '
$ws.Range('E6').Value = '  -1.35%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('D8').Value = '32.74'
$ws.Range('E8').Value = '  +7.07%  '
$ws.Range('E9').Value = '  +0.56%  '
$ws.Range('D10').Value = '0.0677'
$ws.Range('E10').Value = '  +1.61%  '
$ws.Range('D11').Value = '0.0936'
$ws.Range('E11').Value = '  +1.53%  '
$ws.Range('D12').Value = '2.043.17'
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').Value = '11.17'
$ws.Range('E13').Value = '  +11.53%  '
$ws.Range('D14').Value = '1.780.03'
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').Value = '34.573.27'
$ws.Range('E15').Value = '  +2.35%  '
$ws.Range('D16').Value = '0.631'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '4.31'
$ws.Range('E17').Value = '  +2.78%  '
$ws.Range('D18').Value = '68.56'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').Value = '253.30'
$ws.Range('E19').Value = '  +0.44%  '
$ws.Range('D20').Value = '0.0₃0772'
$ws.Range('E20').Value = '  +4.30%  '
$ws.Range('E21').Value = '  +0.04%  '
$ws.Range('D22').Value = '10.44'
$ws.Range('E22').Value = '  +1.43%  '
$ws.Range('E23').Value = '  +0.84%  '
$ws.Range('E24').Value = '  -0.27%  '
$ws.Range('D25').Value = '158.53'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('D26').Value = '16.33'
$ws.Range('E26').Value = '  -1.12%  '
$ws.Range('D27').Value = '7.09'
$ws.Range('E27').Value = '  +1.93%  '
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').Value = '3.75'
$ws.Range('E30').Value = '  -1.97%  '
$ws.Range('E31').Value = '  +0.14%  '
$ws.Range('E32').Value = '  -0.67%  '
$ws.Range('D33').Value = '3.57'
$ws.Range('E33').Value = '  +0.56%  '
$ws.Range('D34').Value = '1.85'
$ws.Range('E34').Value = '  -0.11%  '
$ws.Range('D35').Value = '1.441.80'
$ws.Range('E35').Value = '  -2.84%  '
$ws.Range('E36').Value = '  -1.24%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.631'
$ws.Range('E37').Value = '  -0.67%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.0189'
$ws.Range('E38').Value = '  +2.28%  '
$ws.Range('D39').Value = '83.03'
$ws.Range('E39').Value = '  -0.34%  '
$ws.Range('D40').Value = '2.81'
$ws.Range('E40').Value = '  +4.39%  '
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('D42').Value = '0.903'
$ws.Range('E42').Value = '  +1.93%  '
$ws.Range('D43').Value = '2.07'
$ws.Range('E43').Value = '  -1.04%  '
$ws.Range('D44').Value = '0.0505'
$ws.Range('E44').Value = '  -1.42%  '
$ws.Range('D45').Value = '5.91'
$ws.Range('E45').Value = '  +3.05%  '
$ws.Range('D46').Value = '1.04'
$ws.Range('E46').Value = '  -2.47%  '
$ws.Range('D47').Value = '1.940.85'
$ws.Range('E47').Value = '  +0.60%  '
$ws.Range('D48').Value = '104.79'
$ws.Range('E48').Value = '  +7.33%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '12.01'
$ws.Range('E49').Value = '  +1.92%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').Value = '1.00'
$ws.Range('E50').Value = '  +0.11%  '
$ws.Range('B51').Value = 'BitcoinSV'
$ws.Range('C51').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D51').Value = '49.37'
$ws.Range('E51').Value = '  -2.74%  '
Write-Output "Applied 96 cell updates"
